# Generate Report for Handback
# Updates the "Ready for handoff" rows to "Handed back: in sync with en-US"
# for the 4b732e66-... file across the Overview, zh-cn and de-de sheets,
# refreshes the handback timestamps, and clears the stale "version not
# latest" error detail now that the handback is current.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K3").Value = "2016-08-28 04:47:58"
$ws.Range("P3").Value = ""

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K3").Value = "2016-08-28 04:48:09"
$ws.Range("P3").Value = ""
